$d = $word.ActiveDocument

# The second paragraph contains a Word field (fldChar begin / instrText.../ fldChar end)
# encoding "{ m:'some text'.asParagraph().setAlignment('CENTER') }". We rewrite it
# as plain literal-text runs using "{" and "}" delimiters instead of real field
# characters (TokenIteratorFieldRewriterSplit style), dropping the leading/trailing
# single-space instrText runs.

$targetPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Fields.Count -gt 0) {
        $targetPara = $p
    }
}

$ns = "xmlns:w=`"http://schemas.openxmlformats.org/wordprocessingml/2006/main`""

$runsXml = (
    "<w:r $ns><w:rPr><w:lang w:val=`"en-US`"/></w:rPr><w:t>{</w:t></w:r>" +
    "<w:r $ns><w:rPr><w:lang w:val=`"en-US`"/></w:rPr><w:t>m</w:t></w:r>" +
    "<w:r $ns><w:rPr><w:lang w:val=`"en-US`"/></w:rPr><w:t>:'</w:t></w:r>" +
    "<w:r $ns><w:rPr><w:lang w:val=`"en-US`"/></w:rPr><w:t>some text</w:t></w:r>" +
    "<w:r $ns><w:rPr><w:lang w:val=`"en-US`"/></w:rPr><w:t>'</w:t></w:r>" +
    "<w:r $ns><w:rPr><w:lang w:val=`"en-US`"/></w:rPr><w:t>.asParagraph()</w:t></w:r>" +
    "<w:r $ns><w:rPr><w:lang w:val=`"en-US`"/></w:rPr><w:t>.</w:t></w:r>" +
    "<w:r $ns><w:rPr><w:lang w:val=`"en-US`"/></w:rPr><w:t>set</w:t></w:r>" +
    "<w:r $ns><w:rPr><w:lang w:val=`"en-US`"/></w:rPr><w:t>Alignment</w:t></w:r>" +
    "<w:r $ns><w:rPr><w:lang w:val=`"en-US`"/></w:rPr><w:t>(</w:t></w:r>" +
    "<w:r $ns><w:rPr><w:lang w:val=`"en-US`"/></w:rPr><w:t>'</w:t></w:r>" +
    "<w:r $ns><w:rPr><w:lang w:val=`"en-US`"/></w:rPr><w:t>CENTER</w:t></w:r>" +
    "<w:r $ns><w:rPr><w:lang w:val=`"en-US`"/></w:rPr><w:t>'</w:t></w:r>" +
    "<w:r $ns><w:rPr><w:lang w:val=`"en-US`"/></w:rPr><w:t>)</w:t></w:r>" +
    "<w:r $ns><w:t xml:space=`"preserve`">}</w:t></w:r>"
)

$paraXml = "<w:p $ns>$runsXml</w:p>"

$targetPara.Range.InsertXML($paraXml)
